$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "30.251.55"
$ws.Cells.Item(2, 5).Value = "  -0.20%  "

$ws.Cells.Item(3, 4).Value = "1.925.13"
$ws.Cells.Item(3, 5).Value = "  -0.27%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.000"
$ws.Cells.Item(4, 5).Value = "  +0.06%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "248.34"
$ws.Cells.Item(5, 5).Value = "  -0.11%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.7138"
$ws.Cells.Item(6, 5).Value = "  -1.01%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.000"
$ws.Cells.Item(7, 5).Value = "  -0.01%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3199"
$ws.Cells.Item(8, 5).Value = "  -3.67%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "27.69"
$ws.Cells.Item(9, 5).Value = "  -2.57%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.07051"
$ws.Cells.Item(10, 5).Value = "  +1.72%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.7907"
$ws.Cells.Item(11, 5).Value = "  -1.72%  "

$ws.Cells.Item(12, 5).Value = "  -1.36%  "

$ws.Cells.Item(13, 4).Value = "1.928.69"
$ws.Cells.Item(13, 5).Value = "  -0.12%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "5.384"
$ws.Cells.Item(14, 5).Value = "  -0.62%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "94.75"
$ws.Cells.Item(15, 5).Value = "  -0.04%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "14.64"
$ws.Cells.Item(16, 5).Value = "  +0.72%  "

$ws.Cells.Item(17, 4).Value = "30.286.94"
$ws.Cells.Item(17, 5).Value = "  +0.08%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "257.13"
$ws.Cells.Item(18, 5).Value = "  +1.49%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.000008050"
$ws.Cells.Item(19, 5).Value = "  -3.35%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "5.745"
$ws.Cells.Item(20, 5).Value = "  -1.11%  "

$ws.Cells.Item(21, 4).Value = "2.190.11"
$ws.Cells.Item(21, 5).Value = "  +0.53%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.9999"
$ws.Cells.Item(22, 5).Value = "  +0.00%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.9995"
$ws.Cells.Item(23, 5).Value = "  -0.01%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "6.840"
$ws.Cells.Item(24, 5).Value = "  -0.43%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "9.530"
$ws.Cells.Item(25, 5).Value = "  -2.18%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "165.80"
$ws.Cells.Item(26, 5).Value = "  +3.96%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "19.11"
$ws.Cells.Item(27, 5).Value = "  -0.24%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.253"
$ws.Cells.Item(28, 5).Value = "  -6.80%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.1256"
$ws.Cells.Item(29, 5).Value = "  -6.39%  "

$ws.Cells.Item(30, 5).Value = "  +1.40%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.524"
$ws.Cells.Item(31, 5).Value = "  -1.82%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.386"
$ws.Cells.Item(32, 5).Value = "  -0.41%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "4.111"
$ws.Cells.Item(33, 5).Value = "  -1.98%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.05127"
$ws.Cells.Item(34, 5).Value = "  +0.27%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.267"
$ws.Cells.Item(35, 5).Value = "  +4.01%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.7433"
$ws.Cells.Item(36, 5).Value = "  +0.18%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.765"
$ws.Cells.Item(37, 5).Value = "  +1.07%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.01953"
$ws.Cells.Item(38, 5).Value = "  -1.41%  "

$ws.Cells.Item(39, 5).Value = "  -1.13%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "77.35"
$ws.Cells.Item(40, 5).Value = "  -2.18%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "6.355"
$ws.Cells.Item(41, 5).Value = "  -4.04%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.4496"
$ws.Cells.Item(42, 5).Value = "  +0.57%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.990"
$ws.Cells.Item(43, 5).Value = "  -0.46%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.8448"
$ws.Cells.Item(44, 5).Value = "  +0.92%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.9995"
$ws.Cells.Item(45, 5).Value = "  -0.12%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "100.44"
$ws.Cells.Item(46, 5).Value = "  -1.90%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "9.756"
$ws.Cells.Item(47, 5).Value = "  +0.00%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "7.424"
$ws.Cells.Item(48, 5).Value = "  +1.56%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "36.53"
$ws.Cells.Item(49, 5).Value = "  +0.02%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.06103"
$ws.Cells.Item(50, 5).Value = "  +2.47%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.4196"
$ws.Cells.Item(51, 5).Value = "  +2.59%  "
